$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 header / Row 5 data - new columns J, K
$ws.Range("J4").Value = "SuccessMsgMasterTmplCreation"
$ws.Range("J4").Font.Bold = $true
$ws.Range("J5").Value = "Successfully created 1 records"

$ws.Range("K4").Value = "SuccessMsg"
$ws.Range("K4").Font.Bold = $true
$ws.Range("K5").Value = "1 record(s) successfully entered."

# Row 1 header / Row 2 data - new column F (SuccessMsg / its sample value)
$ws.Range("F1").Value = "SuccessMsg"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F2").Value = "1 record(s) successfully entered."

# Update selection to match new active cell location
$ws.Range("K6").Select()
